$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 21: FEROGLOBIN 30 CAPS ---
$ws.Range("H21").Value = "0:0"
# P21 holds a text value formatted like a number ("270.0000") under a numeric
# display format (0.00). Force the cell to Text so Excel keeps it as a
# shared-string instead of silently coercing it into a numeric value, then
# restore the original numeric display format (cosmetic only; the cached
# value stays text either way, matching the source workbook).
$ws.Range("P21").NumberFormat = "@"
$ws.Range("P21").Value = "270.0000"
$ws.Range("P21").NumberFormat = "0.00"
$ws.Range("Q21").Value = "1:1"

# --- Row 49: was "سائل ريد", now "سرنجات 3 سم" ---
$ws.Range("C49").Value = "سرنجات 3 سم"
$ws.Range("H49").Value = "0:0"
$ws.Range("N49").Value = "2.00"
$ws.Range("P49").NumberFormat = "@"
$ws.Range("P49").Value = "38.0000"
$ws.Range("P49").NumberFormat = "0.00"
$ws.Range("Q49").Value = "19:0"

# --- Row 50: was "سرنجات 3 سم", now "شاش فازلين 10*10 سم" ---
$ws.Range("C50").Value = "شاش فازلين 10*10 سم"
$ws.Range("H50").Value = "11:0"
$ws.Range("N50").Value = "7.00"
$ws.Range("P50").NumberFormat = "@"
$ws.Range("P50").Value = "7.0000"
$ws.Range("P50").NumberFormat = "0.00"
$ws.Range("Q50").Value = "1:0"

# --- Row 51: was "شاش فازلين 10*10 سم", now new item "صابونه دوف SOAP" ---
$ws.Range("C51").Value = "صابونه دوف SOAP"
$ws.Range("H51").Value = "5:0"
$ws.Range("N51").Value = "40.00"
$ws.Range("P51").NumberFormat = "@"
$ws.Range("P51").Value = "40.0000"
$ws.Range("P51").NumberFormat = "0.00"
# Q51 stays "1:0" - unchanged

# --- Row 53: grand total ---
$ws.Range("P53").Value = 3220.4549999999999

# --- Footer timestamp ---
$ws.Range("A54").Value = "Sunday, 1 June, 2025 1:17 PM"
